$d = $word.ActiveDocument

# Step 1: shorten the intro clause "El juego se puede llevar a cabo con dos
# rondas y una tercera en caso de empate. Un jugador gana una ronda cuando"
# down to "Un jugador gana cuando".
$find = $d.Content.Find
$find.Execute(
    "El juego se puede llevar a cabo con dos rondas y una tercera en caso de empate. Un jugador gana una ronda cuando",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Un jugador gana cuando", 2)

# Step 2: move the paragraph break up: remove the trailing
# " y el jugador que gane dos rondas gana la partida." clause (together with
# the original paragraph mark that followed it) and replace it with a period
# followed by a new paragraph mark, so the sentence now ends right after
# "cuadros de vida" and a fresh paragraph begins with "Ya que es un juego...".
$full = $d.Content.Text
$startMarker = " y el jugador que gane dos rondas gana la partida."
$endMarker = "Ya que es un juego"
$start = $full.IndexOf($startMarker)
$end = $full.IndexOf($endMarker)
$r = $d.Range($start, $end)
$r.Text = ".`r"
